$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H (gender), shifting gender/dob/religion right
# so gender->I, dob->J, religion->K
$ws.Range("H1").EntireColumn.Insert()

# Copy formatting from the neighboring header cell (G1) onto the new header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Fotoana"

# Update row 2 values
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0346117893"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "RABEMILSON FRED"
$ws.Range("G2").Value = "Miasa"
$ws.Range("H2").Value = "Jeudi "
$ws.Range("K2").Value = "hindouisme"

# Delete rows 3 and 4 (old duplicate data rows)
$ws.Range("A3:A4").EntireRow.Delete()

$excel.CutCopyMode = $false
